$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $range = $ws.Range($cell)
    $origStyle = $range.Style
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = $origStyle
}

Set-TextValue "D2" '26.996.43'
Set-TextValue "E2" '  -0.06%  '
Set-TextValue "D3" '1.564.29'
Set-TextValue "E3" '  +0.36%  '
Set-TextValue "E4" '  +0.14%  '
Set-TextValue "D5" '207.77'
Set-TextValue "E5" '  +0.02%  '
Set-TextValue "E6" '  -0.01%  '
Set-TextValue "E7" '  +0.05%  '
Set-TextValue "E8" '  +0.19%  '
Set-TextValue "E9" '  +0.00%  '
Set-TextValue "E10" '  +2.10%  '
Set-TextValue "E11" '  -0.49%  '
Set-TextValue "D12" '1.786.18'
Set-TextValue "E12" '  +0.12%  '
Set-TextValue "D13" '1.561.63'
Set-TextValue "E13" '  -0.07%  '
Set-TextValue "D14" '3.75'
Set-TextValue "E14" '  +0.03%  '
Set-TextValue "E15" '  +0.06%  '
Set-TextValue "D16" '62.02'
Set-TextValue "D17" '26.991.24'
Set-TextValue "E17" '  -0.18%  '
Set-TextValue "D18" '0.0₃0703'
Set-TextValue "E18" '  +1.01%  '
Set-TextValue "D19" '216.12'
Set-TextValue "E19" '  -1.49%  '
Set-TextValue "D20" '7.35'
Set-TextValue "E20" '  +0.32%  '
Set-TextValue "E21" '  +0.15%  '
Set-TextValue "E22" '  +1.39%  '
Set-TextValue "D23" '9.19'
Set-TextValue "E23" '  -0.53%  '
Set-TextValue "E24" '  -0.39%  '
Set-TextValue "E25" '  -0.87%  '
Set-TextValue "E26" '  +0.22%  '
Set-TextValue "D27" '15.09'
Set-TextValue "E27" '  +0.79%  '
Set-TextValue "E28" '  +1.39%  '
Set-TextValue "E29" '  +0.09%  '
Set-TextValue "E30" '  +0.39%  '
Set-TextValue "E31" '  +1.20%  '
Set-TextValue "E32" '  +0.05%  '
Set-TextValue "E33" '  +1.32%  '
Set-TextValue "D34" '1.424.27'
Set-TextValue "E34" '  -1.40%  '
Set-TextValue "E35" '  +2.88%  '
Set-TextValue "E36" '  +10.37%  '
Set-TextValue "E37" '  +2.00%  '
Set-TextValue "E38" '  -0.31%  '
Set-TextValue "D39" '0.532'
Set-TextValue "E39" '  +1.59%  '
Set-TextValue "D40" '5.82'
Set-TextValue "E40" '  +1.95%  '
Set-TextValue "E41" '  -0.78%  '
Set-TextValue "E42" '  +0.15%  '
Set-TextValue "E43" '  +1.90%  '
Set-TextValue "E44" '  +1.84%  '
Set-TextValue "D45" '64.78'
Set-TextValue "E45" '  +0.83%  '
Set-TextValue "E46" '  -1.61%  '
Set-TextValue "D47" '1.700.23'
Set-TextValue "E47" '  +0.09%  '
Set-TextValue "D48" '87.35'
Set-TextValue "E48" '  +0.72%  '
Set-TextValue "D49" '0.0519'
Set-TextValue "E49" '  -1.05%  '
Set-TextValue "E50" '  +1.72%  '
Set-TextValue "D51" '0.0959'
Set-TextValue "E51" '  -0.53%  '
